## Word COM-interop script applying the tracked changes to report.docx
##
## Summary of edits:
##  1. Move the "_GoBack" bookmark from the very end of the document (right
##     after the last inline image) to the empty paragraph that follows the
##     "...total approximation errors." paragraph near the top of the report.
##     (Adding a bookmark with an existing name re-seats it, so a single
##     Bookmarks.Add call both creates the new one and removes the old one.)
##  2. Extend the "Demo" heading to read "Demo and Source".
##  3. After the existing "...using the link." sentence, add a new sentence
##     "The code is published in my GitHub repo" where "GitHub repo" is a
##     new hyperlink pointing at the project's GitHub repository.
##
## NOTE: position arithmetic (Range.Start / Range.End) gets stale the moment
## the document is mutated, so every step below re-locates its target with
## Find against the live document instead of reusing old offsets.

$d = $word.ActiveDocument

## ------------------------------------------------------------------
## 1. Re-seat the "_GoBack" bookmark onto the paragraph right after
##    "...total approximation errors."
## ------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("total approximation errors.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng.Collapse(0)
$targetPara = $rng.Paragraphs(1).Next()
$bmRange = $targetPara.Range
$bmRange.Collapse(1)
$d.Bookmarks.Add("_GoBack", $bmRange)

## ------------------------------------------------------------------
## 2. "Demo" -> "Demo and Source"
## ------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("Demo", $true, $true, $false, $false, $false, $true, 1, $false, "", 0)
$rng.Collapse(0)
$rng.InsertAfter(" and Source")

## ------------------------------------------------------------------
## 3. Append the GitHub-repo sentence and hyperlink after the existing
##    "the link" hyperlink / trailing period.
## ------------------------------------------------------------------
$existingLink = $d.Hyperlinks(1)
$afterLink = $existingLink.Range.End
# step past the trailing "." run that already follows the hyperlink
$insertPoint = $d.Range($afterLink + 1, $afterLink + 1)
$insertPoint.InsertAfter(" The code is published in my GitHub repo")

$rng = $d.Content
$rng.Find.Execute("GitHub repo", $true, $true, $false, $false, $false, $true, 1, $false, "", 0)
$d.Hyperlinks.Add($rng, "https://github.com/vanishmax/solving-de")
